$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Convert the old "estimated rf all" average formulas (H10:J10) into the
#     literal (recomputed) values that come from appending the new GTset ---
$ws.Range("H10").Value = 0.91590214067278297
$ws.Range("I10").Value = 0.36939070930385598
$ws.Range("J10").Value = 1.7071606697541899

# --- Add the new "Combined" section (header + per-file rows + aggregate row) ---
$ws.Range("A13").Value = "Combined"

$ws.Range("A14").Value = 1
$ws.Range("B14").Value = "GS"
$ws.Range("C14").Value = "rf AW15_AU_BS3_files_705-749"
$ws.Range("D14").Value = 1000
$ws.Range("E14").Value = 639
$ws.Range("F14").Value = 361
$ws.Range("G14").Value = 58
$ws.Range("H14").Value = 0.91678622668579601
$ws.Range("I14").Value = 0.36099999999999999
$ws.Range("J14").Value = 1.7700831024930701
$ws.Range("K14").Value = 0.94777387726579698

$ws.Range("A15").Value = 1
$ws.Range("B15").Value = "GS"
$ws.Range("C15").Value = "rf BS12_AU_02a_files_1-46"
$ws.Range("D15").Value = 516
$ws.Range("E15").Value = 393
$ws.Range("F15").Value = 123
$ws.Range("G15").Value = 8
$ws.Range("H15").Value = 0.98004987531172105
$ws.Range("I15").Value = 0.23837209302325599
$ws.Range("J15").Value = 3.1951219512195101
$ws.Range("K15").Value = 0.94777387726579698

$ws.Range("A16").Value = 1
$ws.Range("B16").Value = "GS"
$ws.Range("C16").Value = "rf AW14_AU_BS3_files_1-71"
$ws.Range("D16").Value = 1489
$ws.Range("E16").Value = 952
$ws.Range("F16").Value = 537
$ws.Range("G16").Value = 67
$ws.Range("H16").Value = 0.93424926398429797
$ws.Range("I16").Value = 0.36064472800537301
$ws.Range("J16").Value = 1.77281191806331
$ws.Range("K16").Value = 0.94777387726579698

$ws.Range("A17").Value = 1
$ws.Range("B17").Value = "GS"
$ws.Range("C17").Value = "rf BS13_AU_04_files_137-224"
$ws.Range("D17").Value = 1416
$ws.Range("E17").Value = 670
$ws.Range("F17").Value = 746
$ws.Range("G17").Value = 19
$ws.Range("H17").Value = 0.97242380261248196
$ws.Range("I17").Value = 0.52683615819208995
$ws.Range("J17").Value = 0.898123324396783
$ws.Range("K17").Value = 0.94777387726579698

$ws.Range("A18").Value = 1
$ws.Range("B18").Value = "GS"
$ws.Range("C18").Value = "rf AW12_AU_BS3_files_1-250"
$ws.Range("D18").Value = 2736
$ws.Range("E18").Value = 1631
$ws.Range("F18").Value = 1105
$ws.Range("G18").Value = 176
$ws.Range("H18").Value = 0.90260099612617595
$ws.Range("I18").Value = 0.403874269005848
$ws.Range("J18").Value = 1.4760180995475101
$ws.Range("K18").Value = 0.94777387726579698

$ws.Range("A19").Value = 1
$ws.Range("B19").Value = "GS"
$ws.Range("C19").Value = "rf AW12_AU_BS3_files_1464-1507"
$ws.Range("D19").Value = 739
$ws.Range("E19").Value = 554
$ws.Range("F19").Value = 185
$ws.Range("G19").Value = 65
$ws.Range("H19").Value = 0.89499192245557302
$ws.Range("I19").Value = 0.25033829499323401
$ws.Range("J19").Value = 2.9945945945945902
$ws.Range("K19").Value = 0.94777387726579698

$ws.Range("A20").Value = 1
$ws.Range("B20").Value = "GS"
$ws.Range("C20").Value = "rf AW14_AU_BS3_files_309-369"
$ws.Range("D20").Value = 859
$ws.Range("E20").Value = 614
$ws.Range("F20").Value = 245
$ws.Range("G20").Value = 104
$ws.Range("H20").Value = 0.85515320334261802
$ws.Range("I20").Value = 0.28521536670547099
$ws.Range("J20").Value = 2.5061224489795899
$ws.Range("K20").Value = 0.94777387726579698

$ws.Range("A21").Value = 1
$ws.Range("B21").Value = "GS"
$ws.Range("C21").Value = "rf all"
$ws.Range("D21").Value = 8755
$ws.Range("E21").Value = 5453
$ws.Range("F21").Value = 3302
$ws.Range("G21").Value = 497
$ws.Range("H21").Value = 0.91647058823529404
$ws.Range("I21").Value = 0.377155910908053
$ws.Range("J21").Value = 1.6514233797698401
$ws.Range("K21").Value = 0.94777387726579698

# --- Update view: selection moves to M20, matching where the new data ends ---
$ws.Range("M20").Select()

# --- Set page orientation (adds <pageSetup orientation="portrait".../>) ---
$ws.PageSetup.Orientation = 1
